$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$s.Delete()
